# Applies the "output generated at 456a3b4" refresh: bumps the
# "想去人数" (want-to-go count) column F for the rows that changed
# on 展览, 演出, and 全部类型 (the aggregate sheet mirrors the per-category rows).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 468
$ws.Range("F6").Value = 222
$ws.Range("F7").Value = 202
$ws.Range("F8").Value = 232
$ws.Range("F9").Value = 2818
$ws.Range("F10").Value = 55
$ws.Range("F11").Value = 116
$ws.Range("F12").Value = 2157
$ws.Range("F13").Value = 240
$ws.Range("F17").Value = 2506
$ws.Range("F19").Value = 1242
$ws.Range("F20").Value = 4513
$ws.Range("F22").Value = 4226
$ws.Range("F23").Value = 1270
$ws.Range("F24").Value = 2762
$ws.Range("F25").Value = 3169
$ws.Range("F26").Value = 129
$ws.Range("F27").Value = 1475
$ws.Range("F28").Value = 231
$ws.Range("F30").Value = 82
$ws.Range("F31").Value = 233
$ws.Range("F32").Value = 827
$ws.Range("F33").Value = 1422
$ws.Range("F34").Value = 107
$ws.Range("F35").Value = 222
$ws.Range("F36").Value = 561
$ws.Range("F37").Value = 148
$ws.Range("F38").Value = 271
$ws.Range("F39").Value = 345

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 88

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 468
$ws.Range("F5").Value = 88
$ws.Range("F8").Value = 222
$ws.Range("F9").Value = 202
$ws.Range("F11").Value = 232
$ws.Range("F12").Value = 2818
$ws.Range("F13").Value = 55
$ws.Range("F14").Value = 116
$ws.Range("F16").Value = 240
$ws.Range("F21").Value = 2506
$ws.Range("F22").Value = 1242
$ws.Range("F26").Value = 4513
$ws.Range("F28").Value = 4226
$ws.Range("F29").Value = 1270
$ws.Range("F30").Value = 2762
$ws.Range("F31").Value = 3169
$ws.Range("F32").Value = 129
$ws.Range("F35").Value = 1475
$ws.Range("F37").Value = 231
$ws.Range("F39").Value = 82
$ws.Range("F40").Value = 233
$ws.Range("F41").Value = 827
$ws.Range("F43").Value = 1422
$ws.Range("F44").Value = 107
$ws.Range("F45").Value = 222
$ws.Range("F46").Value = 561
$ws.Range("F47").Value = 148
$ws.Range("F48").Value = 271
$ws.Range("F49").Value = 345
